$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates (C1, D1, E1) ---
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# --- Data rows: column C becomes the species string (was numeric duplicate of B) ---
# Column D stays the same species string (unchanged)
# Column E becomes a numeric rejection-f score (was the species string)

$speciesName = "s__CAG-988 sp003149915"

$ws.Range("C2").Value = $speciesName
$ws.Range("C3").Value = $speciesName
$ws.Range("C4").Value = $speciesName
$ws.Range("C5").Value = $speciesName
$ws.Range("C6").Value = $speciesName
$ws.Range("C7").Value = $speciesName
$ws.Range("C8").Value = $speciesName
$ws.Range("C9").Value = $speciesName

$ws.Range("E2").Value = 0.9989549373115363
$ws.Range("E3").Value = 0.9991192576772857
$ws.Range("E4").Value = 0.9990260119264609
$ws.Range("E5").Value = 0.9990242227992739
$ws.Range("E6").Value = 0.9990618712264822
$ws.Range("E7").Value = 0.9991198245047196
$ws.Range("E8").Value = 0.9991240693994353
$ws.Range("E9").Value = 0.9989773265358727
